$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update harvester column (B) for data rows 2-14
$ws.Range("B2:B14").Value = "S.GISH"

# Add experimentDesign column (D) for data rows 2-14
$ws.Range("D2:D14").Value = "90minuteInduction"

# Update selection to match the authored view state
$ws.Range("D3:D14").Select()
